$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The shared strings table was reordered (reran text processing),
# which shifts which word appears in column A for a block of rows
# while the counts in column B (and the row order itself) stay fixed.
# Apply the resulting new word for each affected row.
$ws.Range("A16").Value = "полотно"
$ws.Range("A17").Value = "колеса"
$ws.Range("A18").Value = "Крымскую соль"
$ws.Range("A19").Value = "говядина"
$ws.Range("A20").Value = "парча"
$ws.Range("A21").Value = "сено"
$ws.Range("A22").Value = "позумент"
$ws.Range("A23").Value = "табак"
$ws.Range("A30").Value = "ладан"
$ws.Range("A31").Value = "сапог"
$ws.Range("A32").Value = "коса"
$ws.Range("A33").Value = "сани"
$ws.Range("A34").Value = "китайка"
$ws.Range("A35").Value = "овца"
$ws.Range("A36").Value = "обод"
$ws.Range("A37").Value = "ром"
$ws.Range("A38").Value = "конь"
$ws.Range("A39").Value = "замок"
$ws.Range("A40").Value = "гвоздь"
$ws.Range("A41").Value = "веревка"
$ws.Range("A42").Value = "горшок"
$ws.Range("A43").Value = "рогожа"
$ws.Range("A44").Value = "платок"
$ws.Range("A45").Value = "скотский кожа"
$ws.Range("A46").Value = "дуга"
$ws.Range("A47").Value = "бечева"
$ws.Range("A48").Value = "нитка"
$ws.Range("A49").Value = "хомут"
$ws.Range("A50").Value = "сковорода"
$ws.Range("A51").Value = "покроми"
$ws.Range("A53").Value = "роза"
$ws.Range("A54").Value = "котел"
$ws.Range("A55").Value = "гумми"
$ws.Range("A56").Value = "брусья"
